$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").Value = 0.3486084367838913
$ws.Range("J3").Value = 0.5738014500075432
$ws.Range("K3").Value = 0.4457564652853099
$ws.Range("L3").Value = 2.653079583015679
